$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "36.552.49", "0.691") are preserved exactly as text, matching the
# source data which is stored as inline strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "36.552.49"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").Value = "2.062.62"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "0.691"
$ws.Range("E5").Value = "  +3.65%  "

# Row 6
$ws.Range("D6").Value = "243.24"
$ws.Range("E6").Value = "  -2.10%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("D8").Value = "52.85"
$ws.Range("E8").Value = "  -7.74%  "

# Row 9
$ws.Range("D9").Value = "58.64"
$ws.Range("E9").Value = "  -2.69%  "

# Row 10
$ws.Range("D10").Value = "0.362"
$ws.Range("E10").Value = "  -6.39%  "

# Row 11
$ws.Range("D11").Value = "0.0753"
$ws.Range("E11").Value = "  -4.30%  "

# Row 12
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").Value = "0.878"
$ws.Range("E13").Value = "  -3.64%  "

# Row 14
$ws.Range("D14").Value = "14.61"
$ws.Range("E14").Value = "  -10.35%  "

# Row 15
$ws.Range("D15").Value = "2.368.78"
$ws.Range("E15").Value = "  +0.58%  "

# Row 16
$ws.Range("E16").Value = "  -6.28%  "

# Row 17
$ws.Range("D17").Value = "2.111.79"
$ws.Range("E17").Value = "  +2.73%  "

# Row 18
$ws.Range("D18").Value = "36.457.90"
$ws.Range("E18").Value = "  -1.99%  "

# Row 19
$ws.Range("D19").Value = "16.57"
$ws.Range("E19").Value = "  -11.59%  "

# Row 20
$ws.Range("D20").Value = "72.20"
$ws.Range("E20").Value = "  -3.59%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0865"
$ws.Range("E21").Value = "  -4.17%  "

# Row 22
$ws.Range("D22").Value = "5.34"
$ws.Range("E22").Value = "  -2.78%  "

# Row 23
$ws.Range("D23").Value = "237.29"
$ws.Range("E23").Value = "  -0.24%  "

# Row 24
$ws.Range("D24").Value = "0.997"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25
$ws.Range("E25").Value = "  -4.73%  "

# Row 26
$ws.Range("D26").Value = "9.50"
$ws.Range("E26").Value = "  -1.93%  "

# Row 27
$ws.Range("D27").Value = "2.13"
$ws.Range("E27").Value = "  -2.48%  "

# Row 28
$ws.Range("D28").Value = "165.21"
$ws.Range("E28").Value = "  -3.00%  "

# Row 29
$ws.Range("D29").Value = "20.46"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  +2.00%  "

# Row 31
$ws.Range("D31").Value = "5.16"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32
$ws.Range("E32").Value = "  -2.66%  "

# Row 33
$ws.Range("D33").Value = "4.65"
$ws.Range("E33").Value = "  +0.49%  "

# Row 34
$ws.Range("D34").Value = "0.0598"
$ws.Range("E34").Value = "  -4.41%  "

# Row 35
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +3.33%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("E37").Value = "  +2.48%  "

# Row 38
$ws.Range("D38").Value = "0.0809"
$ws.Range("E38").Value = "  -8.56%  "

# Row 39
$ws.Range("D39").Value = "1.25"
$ws.Range("E39").Value = "  -7.08%  "

# Row 40
$ws.Range("D40").Value = "4.82"
$ws.Range("E40").Value = "  -9.30%  "

# Row 41
$ws.Range("D41").Value = "0.0217"
$ws.Range("E41").Value = "  -3.03%  "

# Row 42
$ws.Range("E42").Value = "  -2.14%  "

# Row 43
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  -6.75%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0939"
$ws.Range("E44").Value = "  -6.42%  "

# Row 45
$ws.Range("D45").Value = "94.44"
$ws.Range("E45").Value = "  -2.38%  "

# Row 46
$ws.Range("D46").Value = "1.386.90"
$ws.Range("E46").Value = "  +8.67%  "

# Row 47
$ws.Range("E47").Value = "  +9.92%  "

# Row 48
$ws.Range("D48").Value = "15.59"
$ws.Range("E48").Value = "  -12.34%  "

# Row 49
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  -3.88%  "

# Row 50
$ws.Range("D50").Value = "2.85"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("D51").Value = "2.258.36"
$ws.Range("E51").Value = "  +0.80%  "
